# Hortaliza, Femacal de La Calera - Apio: add a new weekly price record.
# A new row is inserted at row 398 (shifting the existing rows 398-418 down
# to 399-419), and the new row is populated with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 398, pushing everything
# below it (old rows 398-418) down by one row.
$ws.Rows.Item(398).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A398").Value = 3
$ws.Range("B398").Value = "Femacal de La Calera"
$ws.Range("C398").Value = "Coquimbo"
$ws.Range("D398").Value = 44753
$ws.Range("E398").Value = 5
$ws.Range("F398").Value = 100112017
$ws.Range("G398").Value = "Apio"
$ws.Range("H398").Value = "Americana (o)"
$ws.Range("I398").Value = "Primera"
$ws.Range("J398").Value = 250
$ws.Range("K398").Value = 10000
$ws.Range("L398").Value = 11000
$ws.Range("M398").Value = 10480
$ws.Range("N398").Value = "`$/docena de matas"
$ws.Range("O398").Value = "Pan de Azúcar"
$ws.Range("P398").Value = 1747
$ws.Range("Q398").Value = 6
$ws.Range("R398").Value = "Hortaliza"
